$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values per row (B..G = ME, MAE, MSE, RMSE, SE, N)
$data = @{
    2  = @(0.5440079918808879, 0.5440079918808879, 0.3686258812966464, 0.6071456837503223, 0.2797714183633027, 14)
    3  = @(0.371620442828777,  0.3874993246909071, 0.1981274257140365, 0.4451150701942549, 0.2550055127533425, 13)
    4  = @(0.3184154172361096, 0.3634107558758037, 0.1745848588293675, 0.4178335300444035, 0.2825786729951437, 12)
    5  = @(0.3975283727593981, 0.4145492339343189, 0.2091250642104284, 0.4573019398717093, 0.2370777989771771, 11)
    6  = @(0.3560805627268342, 0.3805953001122874, 0.1781696942431756, 0.4221015212519088, 0.2389242722689363, 10)
    7  = @(0.3392710727664942, 0.3699241520168967, 0.1693631219851699, 0.4115375098155329, 0.2470638456253208, 9)
    8  = @(0.352534700000344,  0.3776340216059933, 0.1795390275838225, 0.4237204592462139, 0.251301328242729,  8)
    9  = @(0.323074747095869,  0.3439312876808158, 0.1492574881197667, 0.3863385666999435, 0.2288235751843726, 7)
    10 = @(0.3581453753460085, 0.3775233390900978, 0.1783793918318441, 0.4223498453081807, 0.2452214067738235, 6)
    11 = @(0.3599616089605324, 0.3766256120414225, 0.1882494404226522, 0.4338772181420132, 0.2708253138491602, 5)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    $ws.Cells.Item($row, 2).Value = $values[0]
    $ws.Cells.Item($row, 3).Value = $values[1]
    $ws.Cells.Item($row, 4).Value = $values[2]
    $ws.Cells.Item($row, 5).Value = $values[3]
    $ws.Cells.Item($row, 6).Value = $values[4]
    $ws.Cells.Item($row, 7).Value = $values[5]
}
